# Update financial figures on the FVCB sheet with newly reported values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement section
$ws.Range("D8").Value = 40300     # Total Revenue
$ws.Range("E8").Value = 32600

$ws.Range("D14").Value = 0        # Non Recurring
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0

$ws.Range("D15").Value = 0        # Others
$ws.Range("E15").Value = 0

$ws.Range("D17").Value = 9400     # Total Operating Expenses
$ws.Range("E17").Value = 6900

$ws.Range("D18").Value = 30900    # Operating Income or Loss
$ws.Range("E18").Value = 25700

$ws.Range("D20").Value = -16400   # Total Other Income/Expenses Net
$ws.Range("E20").Value = -15200

$ws.Range("D21").Value = 15100    # Earnings Before Interest And Taxes
$ws.Range("E21").Value = 11000

$ws.Range("D22").Value = 0        # Interest Expense
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0

$ws.Range("D23").Value = 14500    # Income Before Tax
$ws.Range("E23").Value = 10500

$ws.Range("D24").Value = 4800     # Income Tax Expense
$ws.Range("E24").Value = 3600

$ws.Range("D26").Value = 9700     # Income After Tax
$ws.Range("E26").Value = 6900

$ws.Range("D27").Value = 9700     # Net Income From Continuing Ops
$ws.Range("E27").Value = 6900

$ws.Range("D29").Value = -2000    # Discontinued Operations

$ws.Range("D32").Value = 16400    # Other Items
$ws.Range("E32").Value = 15200

$ws.Range("D33").Value = 7700     # Net Income
$ws.Range("E33").Value = 6900

$ws.Range("D35").Value = 7700     # Net Income Applicable To Common Shares
$ws.Range("E35").Value = 6900

# Cash Flow Statement section
$ws.Range("D81").Value = 7700     # Net Income
$ws.Range("E81").Value = 6900

$ws.Range("J101").Value = "NA"    # Effect Of Exchange Rate Changes
